# Add 4 new library rows (93-96) for the rab-3 WT (spike-ins added) samples,
# as described in the commit:
#   "update design matrxi by addign rab-3 WT and update the
#    scripts/prepare_ExpressionMatrix.R by adding rab-3 WT"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for the 4 new rows (columns A-I):
#   A: date (Library date)      B: Request       C: Multiplex
#   D: sample.ID                E: genotype      F: Tissue/Cell-type
#   G: promoter                 H: sampleInfo    I: Tissue.Cell-type.details
$newRows = @(
    @{ Row = 93; A = 43320; B = 6585; C = 6596; D = 71820; E = "WT"; F = "Pan.neurons"; G = "rab-3"; H = "rab-3:HEN1 MLC416 No Treatment L1 (prep 9) WT, spike-ins added"; I = "WT, spike-ins added" },
    @{ Row = 94; A = 43320; B = 6585; C = 6596; D = 71821; E = "WT"; F = "Pan.neurons"; G = "rab-3"; H = "rab-3:HEN1 MLC416 Treatment L1 (prep 9) WT, spike-ins added";    I = "WT, spike-ins added" },
    @{ Row = 95; A = 43320; B = 6585; C = 6596; D = 71822; E = "WT"; F = "Pan.neurons"; G = "rab-3"; H = "rab-3:HEN1 MLC417 No Treatment L1 (prep 6) WT, spike-ins added"; I = "WT, spike-ins added" },
    @{ Row = 96; A = 43320; B = 6585; C = 6596; D = 71823; E = "WT"; F = "Pan.neurons"; G = "rab-3"; H = "rab-3:HEN1 MLC417 Treatment L1 (prep 6) WT, spike-ins added";    I = "WT, spike-ins added" }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $cellA = $ws.Cells($row, 1)
    $cellA.Value = $r.A
    $cellA.NumberFormat = "m/d/yy"
    $cellA.HorizontalAlignment = -4108   # xlCenter

    $cellB = $ws.Cells($row, 2)
    $cellB.Value = $r.B
    $cellB.HorizontalAlignment = -4108

    $cellC = $ws.Cells($row, 3)
    $cellC.Value = $r.C
    $cellC.HorizontalAlignment = -4108

    $cellD = $ws.Cells($row, 4)
    $cellD.Value = $r.D
    $cellD.HorizontalAlignment = -4108

    $cellE = $ws.Cells($row, 5)
    $cellE.Value = $r.E
    $cellE.HorizontalAlignment = -4108

    $cellF = $ws.Cells($row, 6)
    $cellF.Value = $r.F
    $cellF.HorizontalAlignment = -4108

    $cellG = $ws.Cells($row, 7)
    $cellG.Value = $r.G
    $cellG.HorizontalAlignment = -4108

    # Set column I (the repeated "WT, spike-ins added" string) before column H
    # of the first new row so the shared-string table grows in the same order
    # as the target workbook (I's string first, then each row's unique H string).
    $cellI = $ws.Cells($row, 9)
    $cellI.Value = $r.I
    $cellI.HorizontalAlignment = -4108

    $cellH = $ws.Cells($row, 8)
    $cellH.Value = $r.H
    $cellH.HorizontalAlignment = -4108
}

# Reflect the final selection left by the edit.
$ws.Range("H96").Select()
